$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the Membrillo/Champion block (rows 54-56),
# pushing the existing rows 54-77 down to 57-80, to add the latest week of data.
$ws.Rows("54:56").Insert()

# Row 54
$ws.Range("A54").Value2 = 5
$ws.Range("B54").Value = 'Macroferia Regional de Talca'
$ws.Range("C54").Value = 'Maule'
$ws.Range("D54").Value2 = 45040
$ws.Range("E54").Value2 = 7
$ws.Range("F54").Value = 'Fruta'
$ws.Range("G54").Value2 = 100104
$ws.Range("H54").Value = 'Frutos de pepita'
$ws.Range("I54").Value2 = 100104003
$ws.Range("J54").Value = 'Membrillo'
$ws.Range("K54").Value = 'Champion'
$ws.Range("L54").Value = 'Especial'
$ws.Range("M54").Value2 = 180
$ws.Range("N54").Value2 = 11000
$ws.Range("O54").Value2 = 11000
$ws.Range("P54").Value2 = 11000
$ws.Range("Q54").Value = '$/caja 18 kilos granel'
$ws.Range("R54").Value = 'Región de O''Higgins'
$ws.Range("S54").Value2 = 611
$ws.Range("T54").Value2 = 18

# Row 55
$ws.Range("A55").Value2 = 5
$ws.Range("B55").Value = 'Macroferia Regional de Talca'
$ws.Range("C55").Value = 'Maule'
$ws.Range("D55").Value2 = 45040
$ws.Range("E55").Value2 = 7
$ws.Range("F55").Value = 'Fruta'
$ws.Range("G55").Value2 = 100104
$ws.Range("H55").Value = 'Frutos de pepita'
$ws.Range("I55").Value2 = 100104003
$ws.Range("J55").Value = 'Membrillo'
$ws.Range("K55").Value = 'Champion'
$ws.Range("L55").Value = 'Primera'
$ws.Range("M55").Value2 = 100
$ws.Range("N55").Value2 = 9000
$ws.Range("O55").Value2 = 9000
$ws.Range("P55").Value2 = 9000
$ws.Range("Q55").Value = '$/caja 18 kilos granel'
$ws.Range("R55").Value = 'Región de O''Higgins'
$ws.Range("S55").Value2 = 500
$ws.Range("T55").Value2 = 18

# Row 56
$ws.Range("A56").Value2 = 5
$ws.Range("B56").Value = 'Macroferia Regional de Talca'
$ws.Range("C56").Value = 'Maule'
$ws.Range("D56").Value2 = 45040
$ws.Range("E56").Value2 = 7
$ws.Range("F56").Value = 'Fruta'
$ws.Range("G56").Value2 = 100104
$ws.Range("H56").Value = 'Frutos de pepita'
$ws.Range("I56").Value2 = 100104003
$ws.Range("J56").Value = 'Membrillo'
$ws.Range("K56").Value = 'Champion'
$ws.Range("L56").Value = 'Segunda'
$ws.Range("M56").Value2 = 200
$ws.Range("N56").Value2 = 7000
$ws.Range("O56").Value2 = 7000
$ws.Range("P56").Value2 = 7000
$ws.Range("Q56").Value = '$/caja 18 kilos granel'
$ws.Range("R56").Value = 'Región de O''Higgins'
$ws.Range("S56").Value2 = 389
$ws.Range("T56").Value2 = 18
